$d = $word.ActiveDocument

$replacements = @(
    ,@("10+55=65", "72-38=34")
    ,@("35-22=13", "84+13=97")
    ,@("43+15=58", "64+35=99")
    ,@("93-4=89", "46+24=70")
    ,@("62+1=63", "33+38=71")
    ,@("3+17=20", "79-19=60")
    ,@("67-10=57", "35-3=32")
    ,@("25+35=60", "12+73=85")
    ,@("10+1=11", "2+87=89")
    ,@("11+86=97", "29+17=46")
    ,@("0+34=34", "71+23=94")
    ,@("94-67=27", "86-43=43")
    ,@("40+47=87", "65-20=45")
    ,@("27-20=7", "38-37=1")
    ,@("6+36=42", "7+20=27")
    ,@("64+7=71", "75-20=55")
    ,@("59-51=8", "46+18=64")
    ,@("84-68=16", "11-7=4")
    ,@("75-27=48", "2+65=67")
    ,@("13-12=1", "18+10=28")
    ,@("51+13=64", "34-6=28")
    ,@("87-66=21", "14+25=39")
    ,@("41+47=88", "21+5=26")
    ,@("39+47=86", "30+16=46")
    ,@("23+7=30", "50+3=53")
    ,@("41+44=85", "20+75=95")
    ,@("0+27=27", "89-76=13")
    ,@("44-32=12", "20+61=81")
    ,@("99-43=56", "67+14=81")
    ,@("26+60=86", "46+7=53")
    ,@("92-44=48", "11+56=67")
    ,@("41+58=99", "30+24=54")
    ,@("96-70=26", "52-33=19")
    ,@("11+74=85", "5+28=33")
    ,@("54-25=29", "56+0=56")
    ,@("33+63=96", "93-12=81")
    ,@("13+36=49", "44+34=78")
    ,@("70+15=85", "19+24=43")
    ,@("80-27=53", "56-44=12")
    ,@("48-5=43", "86+12=98")
    ,@("98-16=82", "90-65=25")
    ,@("98-31=67", "89-72=17")
    ,@("50-41=9", "9+27=36")
    ,@("76-67=9", "72-1=71")
    ,@("63+12=75", "88-77=11")
    ,@("19+28=47", "29-7=22")
    ,@("1+22=23", "40-6=34")
    ,@("88-29=59", "6+63=69")
    ,@("88-67=21", "75-14=61")
    ,@("25-17=8", "76-2=74")
    ,@("40+54=94", "81-48=33")
    ,@("16+29=45", "42+46=88")
    ,@("39+20=59", "52+15=67")
    ,@("2+16=18", "85-61=24")
    ,@("34+59=93", "10+14=24")
    ,@("31+53=84", "60-57=3")
    ,@("5+83=88", "39+18=57")
    ,@("53+5=58", "94+5=99")
    ,@("93+2=95", "16+18=34")
    ,@("70+0=70", "37+40=77")
    ,@("66-57=9", "70-17=53")
    ,@("57+16=73", "66-22=44")
    ,@("36+29=65", "67+31=98")
    ,@("73+7=80", "67+0=67")
    ,@("28+63=91", "13+1=14")
    ,@("11+53=64", "76+19=95")
    ,@("80-41=39", "30-29=1")
    ,@("18-11=7", "19+40=59")
    ,@("52-17=35", "10+80=90")
    ,@("99-7=92", "36-18=18")
    ,@("41+8=49", "65+22=87")
    ,@("63-2=61", "18+3=21")
    ,@("31+59=90", "25+34=59")
    ,@("98-9=89", "67-60=7")
    ,@("43+1=44", "84-58=26")
    ,@("86-52=34", "79-3=76")
    ,@("33-17=16", "97-1=96")
    ,@("86-0=86", "17-17=0")
    ,@("18+80=98", "54+40=94")
    ,@("44+0=44", "45-32=13")
    ,@("81-55=26", "73-39=34")
    ,@("17+1=18", "92-1=91")
    ,@("62-44=18", "1+51=52")
    ,@("36+3=39", "16+67=83")
    ,@("25+53=78", "11+38=49")
    ,@("65+23=88", "82-33=49")
    ,@("65+26=91", "52-5=47")
    ,@("68-6=62", "67+28=95")
    ,@("80-68=12", "14+19=33")
    ,@("56+36=92", "9+80=89")
    ,@("5-2=3", "77-38=39")
    ,@("87-21=66", "98-66=32")
    ,@("21+28=49", "63+25=88")
    ,@("51-3=48", "91-52=39")
    ,@("54-29=25", "81+9=90")
    ,@("61+29=90", "76-61=15")
    ,@("54+4=58", "7+32=39")
    ,@("81+17=98", "96+3=99")
    ,@("91-91=0", "30+55=85")
    ,@("14+78=92", "6+3=9")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
